# Timing Category workbook update
# - Update several description/label cells per the revised response coding
# - Clear out the last two data rows (What If / Hindcast) that no longer apply
# - Shrink the Timing_Category defined name range to match the new data extent

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Historic Lookback -> Recent Past
$ws.Range("B2").Value = "Recent Past"
$ws.Range("C2").Value = "Need to look back at recent data from model for (e.g.) Forensics, drift analysis, water quality."

# Row 3 - Real Time description tweak
$ws.Range("C3").Value = "Data needed about what is going on now, or data needed as it is generated"

# Row 4 - Prediction description tweak
$ws.Range("C4").Value = "Looking forward a few days, including estimated current conditions"

# Row 5 - Extreme Event description tweak (typo fix + example)
$ws.Range("C5").Value = 'Information about behavior under a specific extreme event, like the "100 year storm"'

# Row 6 - Scenario description tweak
$ws.Range("C6").Value = "Scenarios other than storms or floods, such as sea level rise scenarios or low flow conditions"

# Row 7 - Long Term -> Pattern or Risk
$ws.Range("B7").Value = "Pattern or Risk"
$ws.Range("C7").Value = 'Analysis of  risk, probability of outcomes, long-term trends or "averages"'

# Row 8 - Unclear or Unknown -> Unclear
$ws.Range("B8").Value = "Unclear"

# Rows 10 and 11 - remove the "What If" and "Hindcast" entries entirely
$ws.Range("A10:C11").ClearContents()

# Shrink the named range to the new data extent (A1:C9)
$wb.Names.Item("Timing_Category").RefersTo = "='Timing_Category'!`$A`$1:`$C`$9"
